$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "End" date for Sugarcane (row 9) from "31st March" to "31st December"
$ws.Range("F9").Value = "31st December"

# Move the active selection to F10 (matches recorded cursor position in diff)
$ws.Range("F10").Select()
